$wb = $excel.ActiveWorkbook

# --- Fix field types on the "classFields" sheet ---
# The parser previously mis-aligned field names with their types; correct the
# "Field Type" column (D) so it matches the actual field (column B) again.
$fields = $wb.Worksheets.Item("classFields")

$fields.Cells.Item(2, 2).Value = "productCount"
$fields.Cells.Item(2, 4).Value = "int"

$fields.Cells.Item(3, 2).Value = "productId"
$fields.Cells.Item(3, 4).Value = "java.lang.Long"

$fields.Cells.Item(4, 2).Value = "price"
$fields.Cells.Item(4, 4).Value = "int"

$fields.Cells.Item(5, 2).Value = "id"
$fields.Cells.Item(5, 4).Value = "java.lang.Long"

$fields.Cells.Item(6, 2).Value = "status"
$fields.Cells.Item(6, 4).Value = "java.lang.String"

$fields.Cells.Item(7, 2).Value = "source"
$fields.Cells.Item(7, 4).Value = "java.lang.String"

$fields.Cells.Item(8, 2).Value = "customerId"
$fields.Cells.Item(8, 4).Value = "java.lang.Long"

# --- Fix and extend the "methodNumberOfLines" sheet ---
# toString() actually spans 3 lines (was recorded as 2), and the constructors
# of Order were previously missing from the parsed output entirely.
$methods = $wb.Worksheets.Item("methodNumberOfLines")

# "Number of Lines" values are stored as text (shared strings) in this
# workbook, just like every other column, so force text formatting while
# assigning the digits, then drop the format override again so the cell
# keeps using the default style.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $methods.Cells.Item(16, 3) "3"

$methods.Cells.Item(17, 1).Value = "pl.piomin.base.domain.Order"
$methods.Cells.Item(17, 2).Value = "Order()"
Set-TextValue $methods.Cells.Item(17, 3) "2"

$methods.Cells.Item(18, 1).Value = "pl.piomin.base.domain.Order"
$methods.Cells.Item(18, 2).Value = "Order(java.lang.Long, java.lang.Long, java.lang.Long, java.lang.String)"
Set-TextValue $methods.Cells.Item(18, 3) "6"

$methods.Cells.Item(19, 1).Value = "pl.piomin.base.domain.Order"
$methods.Cells.Item(19, 2).Value = "Order(java.lang.Long, java.lang.Long, java.lang.Long, int, int)"
Set-TextValue $methods.Cells.Item(19, 3) "8"
